$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 261.83334
$ws.Range("I9").Value = 166.33333
$ws.Range("K9").Value = 166.33333
$ws.Range("M9").Value = 2.666670000000011
$ws.Range("H17").Value = 341.0263
$ws.Range("J17").Value = 288.19446
$ws.Range("L17").Value = 864.58338
$ws.Range("N17").Value = -1200.58338
$ws.Range("H32").Value = 25003300
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 25003300
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 25003300
$ws.Range("N32").Value = -25003952
$ws.Range("H33").Value = 339.35715
$ws.Range("I33").Value = 286.1
$ws.Range("K33").Value = 286.1
$ws.Range("M33").Value = -57.10000000000002
$ws.Range("H40").Value = 9399.4
$ws.Range("J40").Value = 9399.4
$ws.Range("L40").Value = 9399.4
$ws.Range("N40").Value = -9749.4
$ws.Range("H41").Value = 200
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 200
$ws.Range("N41").Value = -1080
$ws.Range("H42").Value = 1493.2
$ws.Range("I42").Value = 177.875
$ws.Range("K42").Value = 533.625
$ws.Range("M42").Value = -303.625
$ws.Range("H43").Value = 4535.7646
$ws.Range("J43").Value = 4127
$ws.Range("L43").Value = 4127
$ws.Range("N43").Value = -4265
$ws.Range("H44").Value = 44250
$ws.Range("J44").Value = 44250
$ws.Range("L44").Value = 44250
$ws.Range("N44").Value = -45174
$ws.Range("H53").Value = 1338.3334
$ws.Range("I53").Value = 1000
$ws.Range("J53").Value = 1507.5
$ws.Range("K53").Value = 1000
$ws.Range("L53").Value = 1507.5
$ws.Range("M53").Value = -363
$ws.Range("N53").Value = -2781.5
$ws.Range("H70").Value = 4540
$ws.Range("J70").Value = 5210
$ws.Range("L70").Value = 15630
$ws.Range("N70").Value = -16170
$ws.Range("H73").Value = 4540
$ws.Range("J73").Value = 5210
$ws.Range("L73").Value = 15630
$ws.Range("N73").Value = -17502
$ws.Range("H74").Value = 13750
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 13333.333
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 13333.333
$ws.Range("M74").Value = -14064
$ws.Range("N74").Value = -15205.333
$ws.Range("H76").Value = 3451.5
$ws.Range("I76").Value = 3451.5
$ws.Range("K76").Value = 3451.5
$ws.Range("M76").Value = -3136.5
$ws.Range("H77").Value = 13750
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 13333.333
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 66666.66500000001
$ws.Range("M77").Value = -70320
$ws.Range("N77").Value = -76026.66500000001
$ws.Range("H79").Value = 3451.5
$ws.Range("I79").Value = 3451.5
$ws.Range("K79").Value = 3451.5
$ws.Range("M79").Value = -2359.5
$ws.Range("I86").Value = 5966.6665
$ws.Range("J86").Value = 5875
$ws.Range("K86").Value = 5966.6665
$ws.Range("L86").Value = 5875
$ws.Range("M86").Value = -4843.6665
$ws.Range("N86").Value = -8121
$ws.Range("I89").Value = 5966.6665
$ws.Range("J89").Value = 5875
$ws.Range("K89").Value = 29833.3325
$ws.Range("L89").Value = 29375
$ws.Range("M89").Value = -24217.3325
$ws.Range("N89").Value = -40607
$ws.Range("H96").Value = 1500.2858
$ws.Range("I96").Value = 1144.375
$ws.Range("J96").Value = 1974.8334
$ws.Range("K96").Value = 3433.125
$ws.Range("L96").Value = 5924.5002
$ws.Range("M96").Value = -2060.125
$ws.Range("N96").Value = -8670.5002
$ws.Range("H103").Value = 725
$ws.Range("J103").Value = 733.3333
$ws.Range("L103").Value = 2199.9999
$ws.Range("N103").Value = -3371.9999
$ws.Range("H107").Value = 670.625
$ws.Range("I107").Value = 697.3333
$ws.Range("J107").Value = 270
$ws.Range("K107").Value = 697.3333
$ws.Range("L107").Value = 270
$ws.Range("M107").Value = 1222.6667
$ws.Range("N107").Value = -4110
$ws.Range("H113").Value = 14351.083
$ws.Range("I113").Value = 6959.8
$ws.Range("J113").Value = 19630.572
$ws.Range("K113").Value = 6959.8
$ws.Range("L113").Value = 19630.572
$ws.Range("M113").Value = -3705.8
$ws.Range("N113").Value = -26138.572
$ws.Range("H116").Value = 8268.0625
$ws.Range("I116").Value = 6259
$ws.Range("J116").Value = 9181.272
$ws.Range("K116").Value = 6259
$ws.Range("L116").Value = 9181.272
$ws.Range("M116").Value = -2817
$ws.Range("N116").Value = -16065.272
$ws.Range("H125").Value = 1198.3334
$ws.Range("I125").Value = 1058
$ws.Range("K125").Value = 9522
$ws.Range("M125").Value = -7062
$ws.Range("H132").Value = 12747.622
$ws.Range("I132").Value = 1254.7073
$ws.Range("K132").Value = 3764.1219
$ws.Range("M132").Value = -1234.1219
$ws.Range("H137").Value = 4869.364
$ws.Range("I137").Value = 4569
$ws.Range("J137").Value = 4957.706
$ws.Range("K137").Value = 13707
$ws.Range("L137").Value = 14873.118
$ws.Range("M137").Value = -11157
$ws.Range("N137").Value = -19973.118
$ws.Range("H138").Value = 2547.6667
$ws.Range("J138").Value = 4239.2856
$ws.Range("L138").Value = 12717.8568
$ws.Range("N138").Value = -22997.8568
$ws.Range("H140").Value = 39833.332
$ws.Range("J140").Value = 39833.332
$ws.Range("L140").Value = 39833.332
$ws.Range("N140").Value = -50193.332
$ws.Range("H141").Value = 8372.941
$ws.Range("I141").Value = 7959
$ws.Range("K141").Value = 23877
$ws.Range("M141").Value = -18697
$ws.Range("M32").ClearContents()
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 40202.445
$ws.Range("I2").Value = 67202.2
$ws.Range("K2").Value = 67202.2
$ws.Range("M2").Value = -67089.2
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 200
$ws.Range("N8").Value = -488
$ws.Range("H12").Value = 348.17648
$ws.Range("I12").Value = 219.84616
$ws.Range("J12").Value = 765.25
$ws.Range("K12").Value = 219.84616
$ws.Range("L12").Value = 765.25
$ws.Range("M12").Value = -46.84616
$ws.Range("N12").Value = -1111.25
$ws.Range("H32").Value = 2292.5435
$ws.Range("I32").Value = 2392.85
$ws.Range("K32").Value = 2392.85
$ws.Range("M32").Value = -2105.85
$ws.Range("H45").Value = 9234.571
$ws.Range("I45").Value = 6556
$ws.Range("K45").Value = 6556
$ws.Range("M45").Value = -6179
$ws.Range("H61").Value = 16114.546
$ws.Range("I61").Value = 5395
$ws.Range("K61").Value = 5395
$ws.Range("M61").Value = -5183
$ws.Range("H74").Value = 11998
$ws.Range("J74").Value = 11998
$ws.Range("L74").Value = 11998
$ws.Range("N74").Value = -13746
$ws.Range("H76").Value = 74982
$ws.Range("J76").Value = 74982
$ws.Range("L76").Value = 74982
$ws.Range("N76").Value = -75658
$ws.Range("H77").Value = 11998
$ws.Range("J77").Value = 11998
$ws.Range("L77").Value = 59990
$ws.Range("N77").Value = -68726
$ws.Range("H79").Value = 74982
$ws.Range("J79").Value = 74982
$ws.Range("L79").Value = 74982
$ws.Range("N79").Value = -77322
$ws.Range("H97").Value = 4476.3335
$ws.Range("I97").Value = 2535.875
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 2535.875
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -2039.875
$ws.Range("N97").Value = -20992
$ws.Range("H98").Value = 35641.168
$ws.Range("J98").Value = 35641.168
$ws.Range("L98").Value = 35641.168
$ws.Range("N98").Value = -41631.168
$ws.Range("H102").Value = 20505
$ws.Range("I102").Value = 10673.333
$ws.Range("J102").Value = 50000
$ws.Range("K102").Value = 10673.333
$ws.Range("L102").Value = 50000
$ws.Range("M102").Value = -9051.333
$ws.Range("N102").Value = -53244
$ws.Range("H110").Value = 5938.25
$ws.Range("I110").Value = 5866.625
$ws.Range("J110").Value = 6081.5
$ws.Range("K110").Value = 5866.625
$ws.Range("L110").Value = 6081.5
$ws.Range("M110").Value = -3821.625
$ws.Range("N110").Value = -10171.5
$ws.Range("H116").Value = 40202.445
$ws.Range("I116").Value = 67202.2
$ws.Range("K116").Value = 67202.2
$ws.Range("M116").Value = -64908.2
$ws.Range("H122").Value = 5103.852
$ws.Range("I122").Value = 4914.7144
$ws.Range("J122").Value = 5307.5386
$ws.Range("K122").Value = 14744.1432
$ws.Range("L122").Value = 15922.6158
$ws.Range("M122").Value = -12294.1432
$ws.Range("N122").Value = -20822.6158
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6470
$ws.Range("H135").Value = 39833.332
$ws.Range("J135").Value = 39833.332
$ws.Range("L135").Value = 39833.332
$ws.Range("N135").Value = -49973.332
$ws.Range("H136").Value = 16114.546
$ws.Range("I136").Value = 5395
$ws.Range("K136").Value = 16185
$ws.Range("M136").Value = -13635
$ws.Range("M8").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 40202.445
$ws.Range("I3").Value = 67202.2
$ws.Range("K3").Value = 67202.2
$ws.Range("M3").Value = -67088.2
$ws.Range("H20").Value = 12786.286
$ws.Range("I20").Value = 7376.75
$ws.Range("K20").Value = 7376.75
$ws.Range("M20").Value = -7129.75
$ws.Range("H105").Value = 1763.875
$ws.Range("I105").Value = 967.2353
$ws.Range("J105").Value = 3698.5715
$ws.Range("K105").Value = 967.2353
$ws.Range("L105").Value = 3698.5715
$ws.Range("M105").Value = 779.7647
$ws.Range("N105").Value = -7192.5715
$ws.Range("H107").Value = 11033.533
$ws.Range("I107").Value = 10507.583
$ws.Range("J107").Value = 13137.333
$ws.Range("K107").Value = 10507.583
$ws.Range("L107").Value = 13137.333
$ws.Range("M107").Value = -8587.583
$ws.Range("N107").Value = -16977.333
$ws.Range("H110").Value = 115000
$ws.Range("J110").Value = 115000
$ws.Range("L110").Value = 115000
$ws.Range("N110").Value = -123180
$ws.Range("H134").Value = 2695.6
$ws.Range("I134").Value = 2848.8572
$ws.Range("J134").Value = 2338
$ws.Range("K134").Value = 8546.5716
$ws.Range("L134").Value = 7014
$ws.Range("M134").Value = -6011.571599999999
$ws.Range("N134").Value = -12084

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3021.9092
$ws.Range("I31").Value = 1850.8334
$ws.Range("J31").Value = 4427.2
$ws.Range("K31").Value = 1850.8334
$ws.Range("L31").Value = 4427.2
$ws.Range("M31").Value = -1555.8334
$ws.Range("N31").Value = -5017.2
$ws.Range("H34").Value = 3021.9092
$ws.Range("I34").Value = 1850.8334
$ws.Range("J34").Value = 4427.2
$ws.Range("K34").Value = 1850.8334
$ws.Range("L34").Value = 4427.2
$ws.Range("M34").Value = -1648.8334
$ws.Range("N34").Value = -4831.2
$ws.Range("H58").Value = 2379.7693
$ws.Range("I58").Value = 1537
$ws.Range("K58").Value = 1537
$ws.Range("M58").Value = -1334
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908
$ws.Range("H105").Value = 2887.25
$ws.Range("I105").Value = 3449.6667
$ws.Range("K105").Value = 3449.6667
$ws.Range("M105").Value = -1702.6667
$ws.Range("H122").Value = 3389.6875
$ws.Range("I122").Value = 2205
$ws.Range("K122").Value = 6615
$ws.Range("M122").Value = -4165
$ws.Range("H132").Value = 2567.9412
$ws.Range("I132").Value = 2282
$ws.Range("K132").Value = 6846
$ws.Range("M132").Value = -4316
$ws.Range("H134").Value = 4623.579
$ws.Range("I134").Value = 4844.353
$ws.Range("J134").Value = 2747
$ws.Range("K134").Value = 14533.059
$ws.Range("L134").Value = 8241
$ws.Range("M134").Value = -11998.059
$ws.Range("N134").Value = -13311
$ws.Range("H136").Value = 2379.7693
$ws.Range("I136").Value = 1537
$ws.Range("K136").Value = 4611
$ws.Range("M136").Value = -2061

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 639861.5
$ws.Range("J37").Value = 639861.5
$ws.Range("L37").Value = 1919584.5
$ws.Range("N37").Value = -1919808.5
$ws.Range("H68").Value = 2414.7144
$ws.Range("J68").Value = 2414.7144
$ws.Range("L68").Value = 7244.1432
$ws.Range("N68").Value = -8866.143199999999
$ws.Range("H69").Value = 19609874
$ws.Range("I69").Value = 1307.7693
$ws.Range("J69").Value = 83337710
$ws.Range("K69").Value = 3923.3079
$ws.Range("L69").Value = 250013130
$ws.Range("M69").Value = -3112.3079
$ws.Range("N69").Value = -250014752
$ws.Range("H71").Value = 2414.7144
$ws.Range("J71").Value = 2414.7144
$ws.Range("L71").Value = 21732.4296
$ws.Range("N71").Value = -29844.4296
$ws.Range("H72").Value = 19609874
$ws.Range("I72").Value = 1307.7693
$ws.Range("J72").Value = 83337710
$ws.Range("K72").Value = 11769.9237
$ws.Range("L72").Value = 750039390
$ws.Range("M72").Value = -7713.923699999999
$ws.Range("N72").Value = -750047502
$ws.Range("H80").Value = 6399.9473
$ws.Range("I80").Value = 7400
$ws.Range("J80").Value = 6042.7856
$ws.Range("K80").Value = 22200
$ws.Range("L80").Value = 18128.3568
$ws.Range("M80").Value = -21264
$ws.Range("N80").Value = -20000.3568
$ws.Range("H83").Value = 6399.9473
$ws.Range("I83").Value = 7400
$ws.Range("J83").Value = 6042.7856
$ws.Range("K83").Value = 66600
$ws.Range("L83").Value = 54385.0704
$ws.Range("M83").Value = -61920
$ws.Range("N83").Value = -63745.0704
$ws.Range("H86").Value = 465
$ws.Range("I86").Value = 460.13333
$ws.Range("J86").Value = 477.16666
$ws.Range("K86").Value = 1380.39999
$ws.Range("L86").Value = 1431.49998
$ws.Range("M86").Value = -194.3999899999999
$ws.Range("N86").Value = -3803.49998
$ws.Range("H89").Value = 465
$ws.Range("I89").Value = 460.13333
$ws.Range("J89").Value = 477.16666
$ws.Range("K89").Value = 4141.19997
$ws.Range("L89").Value = 4294.49994
$ws.Range("M89").Value = 1786.80003
$ws.Range("N89").Value = -16150.49994
$ws.Range("H103").Value = 20835062
$ws.Range("I103").Value = 1060
$ws.Range("J103").Value = 66669868
$ws.Range("K103").Value = 3180
$ws.Range("L103").Value = 200009604
$ws.Range("M103").Value = -2301
$ws.Range("N103").Value = -200011362
$ws.Range("H121").Value = 4255.5625
$ws.Range("I121").Value = 781.6667
$ws.Range("J121").Value = 6339.9
$ws.Range("K121").Value = 2345.0001
$ws.Range("L121").Value = 19019.7
$ws.Range("M121").Value = -1035.0001
$ws.Range("N121").Value = -21639.7
$ws.Range("H122").Value = 4965.5625
$ws.Range("I122").Value = 513
$ws.Range("J122").Value = 5993.077
$ws.Range("K122").Value = 4617
$ws.Range("L122").Value = 53937.693
$ws.Range("M122").Value = -2167
$ws.Range("N122").Value = -58837.693
$ws.Range("H140").Value = 626932.3
$ws.Range("I140").Value = 715994.1
$ws.Range("K140").Value = 2147982.3
$ws.Range("M140").Value = -2142802.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 118074.7
$ws.Range("I70").Value = 192483.33
$ws.Range("K70").Value = 192483.33
$ws.Range("M70").Value = -192213.33
$ws.Range("H73").Value = 118074.7
$ws.Range("I73").Value = 192483.33
$ws.Range("K73").Value = 192483.33
$ws.Range("M73").Value = -191547.33
$ws.Range("H74").Value = 20031
$ws.Range("J74").Value = 20031
$ws.Range("L74").Value = 20031
$ws.Range("N74").Value = -21903
$ws.Range("H77").Value = 20031
$ws.Range("J77").Value = 20031
$ws.Range("L77").Value = 60093
$ws.Range("N77").Value = -69453
$ws.Range("H80").Value = 131322.78
$ws.Range("J80").Value = 10159
$ws.Range("N80").Value = -12155
$ws.Range("H83").Value = 131322.78
$ws.Range("J83").Value = 10159
$ws.Range("L83").Value = 50795
$ws.Range("N83").Value = -60779
$ws.Range("H102").Value = 8250
$ws.Range("I102").Value = 6857.143
$ws.Range("K102").Value = 6857.143
$ws.Range("M102").Value = -5235.143
$ws.Range("H113").Value = 9113.333
$ws.Range("J113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("N113").Value = -24340
$ws.Range("H117").Value = 53325
$ws.Range("J117").Value = 53325
$ws.Range("L117").Value = 53325
$ws.Range("N117").Value = -60209
$ws.Range("H122").Value = 4611.846
$ws.Range("I122").Value = 3995.2104
$ws.Range("J122").Value = 6285.5713
$ws.Range("K122").Value = 11985.6312
$ws.Range("L122").Value = 18856.7139
$ws.Range("M122").Value = -9535.6312
$ws.Range("N122").Value = -23756.7139
$ws.Range("H132").Value = 8236.151
$ws.Range("I132").Value = 7426.8335
$ws.Range("J132").Value = 10394.333
$ws.Range("K132").Value = 22280.5005
$ws.Range("L132").Value = 31182.999
$ws.Range("M132").Value = -19750.5005
$ws.Range("N132").Value = -36242.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1902
$ws.Range("I7").Value = 1902
$ws.Range("K7").Value = 1902
$ws.Range("M7").Value = -1790
$ws.Range("H21").Value = 2586
$ws.Range("I21").Value = 502
$ws.Range("J21").Value = 3975.3333
$ws.Range("K21").Value = 502
$ws.Range("L21").Value = 3975.3333
$ws.Range("M21").Value = -328
$ws.Range("N21").Value = -4323.3333
$ws.Range("H40").Value = 200000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 200000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 200000
$ws.Range("N40").Value = -200272
$ws.Range("H46").Value = 2742
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312
$ws.Range("H57").Value = 23000
$ws.Range("I57").Value = 28000
$ws.Range("J57").Value = 18000
$ws.Range("K57").Value = 28000
$ws.Range("L57").Value = 18000
$ws.Range("M57").Value = -27434
$ws.Range("N57").Value = -19132
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H68").Value = 4710.8
$ws.Range("I68").Value = 4999
$ws.Range("J68").Value = 4638.75
$ws.Range("K68").Value = 4999
$ws.Range("L68").Value = 4638.75
$ws.Range("M68").Value = -4250
$ws.Range("N68").Value = -6136.75
$ws.Range("H71").Value = 4710.8
$ws.Range("I71").Value = 4999
$ws.Range("J71").Value = 4638.75
$ws.Range("K71").Value = 24995
$ws.Range("L71").Value = 23193.75
$ws.Range("M71").Value = -21251
$ws.Range("N71").Value = -30681.75
$ws.Range("H82").Value = 37038844
$ws.Range("I82").Value = 55557520
$ws.Range("J82").Value = 1494.4445
$ws.Range("K82").Value = 55557520
$ws.Range("L82").Value = 1494.4445
$ws.Range("M82").Value = -55557159
$ws.Range("N82").Value = -2216.4445
$ws.Range("H85").Value = 37038844
$ws.Range("I85").Value = 55557520
$ws.Range("J85").Value = 1494.4445
$ws.Range("K85").Value = 55557520
$ws.Range("L85").Value = 1494.4445
$ws.Range("M85").Value = -55556272
$ws.Range("N85").Value = -3990.4445
$ws.Range("H93").Value = 1957.1666
$ws.Range("I93").Value = 1848.6
$ws.Range("K93").Value = 1848.6
$ws.Range("M93").Value = -600.5999999999999
$ws.Range("H110").Value = 34749
$ws.Range("J110").Value = 34749
$ws.Range("L110").Value = 34749
$ws.Range("N110").Value = -42929
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H122").Value = 6389.4
$ws.Range("I122").Value = 6432.8887
$ws.Range("J122").Value = 5998
$ws.Range("K122").Value = 19298.6661
$ws.Range("L122").Value = 17994
$ws.Range("M122").Value = -16848.6661
$ws.Range("N122").Value = -22894
$ws.Range("H126").Value = 1902
$ws.Range("I126").Value = 1902
$ws.Range("K126").Value = 5706
$ws.Range("M126").Value = -3236
$ws.Range("H132").Value = 3035.3333
$ws.Range("I132").Value = 2544.6667
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 7634.000100000001
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -5104.000100000001
$ws.Range("N132").Value = -20054
$ws.Range("H136").Value = 2979.9
$ws.Range("I136").Value = 2755.4443
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8266.332900000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5716.332900000001
$ws.Range("N136").Value = -20100
$ws.Range("M40").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 33165
$ws.Range("J49").Value = 33165
$ws.Range("L49").Value = 33165
$ws.Range("N49").Value = -33625
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H96").Value = 250500
$ws.Range("I96").Value = 250500
$ws.Range("K96").Value = 250500
$ws.Range("M96").Value = -249127
$ws.Range("H100").Value = 732.2857
$ws.Range("I100").Value = 1051
$ws.Range("J100").Value = 493.25
$ws.Range("K100").Value = 2102
$ws.Range("L100").Value = 986.5
$ws.Range("M100").Value = -1561
$ws.Range("N100").Value = -2068.5
$ws.Range("H107").Value = 823
$ws.Range("I107").Value = 837
$ws.Range("J107").Value = 776.3333
$ws.Range("K107").Value = 2511
$ws.Range("L107").Value = 2328.9999
$ws.Range("M107").Value = -591
$ws.Range("N107").Value = -6168.9999
$ws.Range("H126").Value = 3937
$ws.Range("I126").Value = 3583
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 10749
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -8279
$ws.Range("N126").Value = -19937
$ws.Range("H132").Value = 4821.186
$ws.Range("I132").Value = 3712.9312
$ws.Range("K132").Value = 11138.7936
$ws.Range("M132").Value = -8608.7936
$ws.Range("H136").Value = 109177.78
$ws.Range("I136").Value = 121887.5
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 365662.5
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -363112.5
$ws.Range("N136").Value = -27600
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
